$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2019.4
$ws.Range("I15").Value = 2019.4
$ws.Range("K15").Value = 6058.200000000001
$ws.Range("M15").Value = -5889.200000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1003.4167
$ws.Range("I28").Value = 735.5
$ws.Range("K28").Value = 735.5
$ws.Range("M28").Value = -250.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2100.4211
$ws.Range("I112").Value = 1666
$ws.Range("K112").Value = 4998
$ws.Range("M112").Value = -3890

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4804.3105
$ws.Range("I137").Value = 2228.3
$ws.Range("K137").Value = 6684.900000000001
$ws.Range("M137").Value = -4134.900000000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2731.4878
$ws.Range("I138").Value = 1449.5
$ws.Range("K138").Value = 4348.5
$ws.Range("M138").Value = 791.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8935735
$ws.Range("I32").Value = 12503268
$ws.Range("J32").Value = 16900.25
$ws.Range("K32").Value = 12503268
$ws.Range("L32").Value = 16900.25
$ws.Range("M32").Value = -12502981
$ws.Range("N32").Value = -17474.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 25002224
$ws.Range("I45").Value = 38463510
$ws.Range("K45").Value = 38463510
$ws.Range("M45").Value = -38463133

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 35720150
$ws.Range("I61").Value = 33337868
$ws.Range("K61").Value = 33337868
$ws.Range("M61").Value = -33337656

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10009494
$ws.Range("I74").Value = 19233466
$ws.Range("K74").Value = 19233466
$ws.Range("M74").Value = -19232592

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 10009494
$ws.Range("I77").Value = 19233466
$ws.Range("K77").Value = 96167330
$ws.Range("M77").Value = -96162962

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 42488.89
$ws.Range("I132").Value = 48046.316
$ws.Range("K132").Value = 144138.948
$ws.Range("M132").Value = -141608.948

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 35720150
$ws.Range("I136").Value = 33337868
$ws.Range("K136").Value = 100013604
$ws.Range("M136").Value = -100011054

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 2719.1428
$ws.Range("J36").Value = 2976
$ws.Range("L36").Value = 2976
$ws.Range("N36").Value = -4044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3838.5386
$ws.Range("I86").Value = 3911.7778
$ws.Range("J86").Value = 3673.75
$ws.Range("K86").Value = 3911.7778
$ws.Range("L86").Value = 3673.75
$ws.Range("M86").Value = -2788.7778
$ws.Range("N86").Value = -5919.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3838.5386
$ws.Range("I89").Value = 3911.7778
$ws.Range("J89").Value = 3673.75
$ws.Range("K89").Value = 19558.889
$ws.Range("L89").Value = 18368.75
$ws.Range("M89").Value = -13942.889
$ws.Range("N89").Value = -29600.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2264.0908
$ws.Range("I105").Value = 1499.5333
$ws.Range("K105").Value = 1499.5333
$ws.Range("M105").Value = 247.4666999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 410855.84
$ws.Range("I134").Value = 1777.5294
$ws.Range("J134").Value = 1106289
$ws.Range("K134").Value = 5332.5882
$ws.Range("L134").Value = 3318867
$ws.Range("M134").Value = -2797.5882
$ws.Range("N134").Value = -3323937

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 579369.3
$ws.Range("J31").Value = 1025078.2
$ws.Range("L31").Value = 1025078.2
$ws.Range("N31").Value = -1025668.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 579369.3
$ws.Range("J34").Value = 1025078.2
$ws.Range("L34").Value = 1025078.2
$ws.Range("N34").Value = -1025482.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2821.8
$ws.Range("I105").Value = 3024.5
$ws.Range("J105").Value = 2011
$ws.Range("K105").Value = 3024.5
$ws.Range("L105").Value = 2011
$ws.Range("M105").Value = -1277.5
$ws.Range("N105").Value = -5505

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1608.8422
$ws.Range("I122").Value = 1621.8182
$ws.Range("J122").Value = 1591
$ws.Range("K122").Value = 4865.4546
$ws.Range("L122").Value = 4773
$ws.Range("M122").Value = -2415.4546
$ws.Range("N122").Value = -9673

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 93997.25
$ws.Range("J127").Value = 78663
$ws.Range("L127").Value = 78663
$ws.Range("N127").Value = -88583

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 13437.714
$ws.Range("I132").Value = 5755.75
$ws.Range("K132").Value = 17267.25
$ws.Range("M132").Value = -14737.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5297.8
$ws.Range("I134").Value = 2163.6667
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 6491.000100000001
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -3956.000100000001
$ws.Range("N134").Value = -35067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 896
$ws.Range("J23").Value = 592.6667
$ws.Range("L23").Value = 1778.0001
$ws.Range("N23").Value = -2248.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1000
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 1000
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 3000
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -3588

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3976.7058
$ws.Range("I118").Value = 1238.7273
$ws.Range("J118").Value = 8996.333000000001
$ws.Range("K118").Value = 3716.1819
$ws.Range("L118").Value = 26988.999
$ws.Range("M118").Value = -2473.1819
$ws.Range("N118").Value = -29474.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1452.6666
$ws.Range("J121").Value = 1861.6875
$ws.Range("L121").Value = 5585.0625
$ws.Range("N121").Value = -8205.0625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 4151.0835
$ws.Range("I138").Value = 3302.1667
$ws.Range("K138").Value = 9906.500100000001
$ws.Range("M138").Value = -4766.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 176291.56
$ws.Range("I141").Value = 337660.88
$ws.Range("K141").Value = 1012982.64
$ws.Range("M141").Value = -1007802.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2087.5454
$ws.Range("J97").Value = 2138.125
$ws.Range("L97").Value = 2138.125
$ws.Range("N97").Value = -3130.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 101084.75
$ws.Range("J110").Value = 101084.75
$ws.Range("L110").Value = 101084.75
$ws.Range("N110").Value = -109264.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3271.9048
$ws.Range("I113").Value = 2553.0908
$ws.Range("J113").Value = 4062.6
$ws.Range("K113").Value = 2553.0908
$ws.Range("L113").Value = 4062.6
$ws.Range("M113").Value = -383.0907999999999
$ws.Range("N113").Value = -8402.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3556.6667
$ws.Range("I126").Value = 3144.2856
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 9432.856800000001
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -6962.856800000001
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 47622076
$ws.Range("I132").Value = 58826564
$ws.Range("K132").Value = 176479692
$ws.Range("M132").Value = -176477162

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12248.904
$ws.Range("I7").Value = 9777.799999999999
$ws.Range("K7").Value = 9777.799999999999
$ws.Range("M7").Value = -9665.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 55556412
$ws.Range("I55").Value = 76923860
$ws.Range("J55").Value = 1058.8
$ws.Range("K55").Value = 76923860
$ws.Range("L55").Value = 1058.8
$ws.Range("M55").Value = -76923687
$ws.Range("N55").Value = -1404.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 12248.904
$ws.Range("I126").Value = 9777.799999999999
$ws.Range("K126").Value = 29333.4
$ws.Range("M126").Value = -26863.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 49995
$ws.Range("J127").Value = 49995
$ws.Range("L127").Value = 49995
$ws.Range("N127").Value = -59915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26457.916
$ws.Range("I54").Value = 18900
$ws.Range("K54").Value = 18900
$ws.Range("M54").Value = -18380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7932.4443
$ws.Range("I81").Value = 3899.5557
$ws.Range("K81").Value = 7799.1114
$ws.Range("M81").Value = -6738.1114

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 7932.4443
$ws.Range("I84").Value = 3899.5557
$ws.Range("K84").Value = 38995.557
$ws.Range("M84").Value = -33691.557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2600
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 326768.88
$ws.Range("I132").Value = 3486.4075
$ws.Range("K132").Value = 10459.2225
$ws.Range("M132").Value = -7929.2225
